## "Entrance added to report"
## Insert a new "Место" (Entrance/Location) column into the header row of
## the report table (row 10), between "Событие" (D10) and
## "Отработано за период" (old E10). The three existing header cells
## E10:G10 shift one column to the right (-> F10:H10) and the new column
## gets the same header formatting as its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the trailing header cells one column to the right, starting from
# the rightmost so we never clobber a value before it has been copied.
$ws.Range("G10").Copy($ws.Range("H10"))
$ws.Range("F10").Copy($ws.Range("G10"))
$ws.Range("E10").Copy($ws.Range("F10"))

# New column: reuse the header formatting (style) of the cell to its left,
# then set its own text.
$ws.Range("D10").Copy($ws.Range("E10"))
$ws.Range("E10").Value = "Место"

# Give the newly used column H the same width as the other data columns.
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Match the saved selection/active cell of the edited workbook.
$ws.Range("H10").Select() | Out-Null
